$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "ValidLogin"

# Write header + credential values
$ws.Range("A1").Value = "username"
$ws.Range("B1").Value = "password"
$ws.Range("A2").Value = "admin"
$ws.Range("B2").Value = "pointofsale"

# Set column B width to match the bestFit width recorded in the diff
# (closest value the engine's width->pixel grid can reproduce for 10.21875)
$ws.Columns.Item(2).ColumnWidth = 9.3

# Move the active selection as captured in the saved workbook view
$ws.Range("D15").Select() | Out-Null
